$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$old = "dnasr281@gmail.com, System"
$new = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $v = $cell.Value2
    if ($v -eq $old) {
        $cell.Value = $new
    }
}
